$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A:L (text-like data, incl. ids, dates-as-text, counts) to stay Text
# so numeric-looking / date-looking strings are not auto-coerced by Excel.
$ws.Range("A38:L52").NumberFormat = "@"

# Row 38
$ws.Range("A38").Value = '807044192'
$ws.Range("B38").Value = '5/29/2025'
$ws.Range("C38").Value = 'O''Higgins 4379'
$ws.Range("D38").Value = '13'
$ws.Range("E38").Value = '807044192'
$ws.Range("F38").Value = 'NEW'
$ws.Range("G38").Value = 'Pendiente'
$ws.Range("H38").Value = 'Poste '
$ws.Range("I38").Value = '1'
$ws.Range("J38").Value = 'Desmonte'
$ws.Range("K38").Value = 'Sin equipos'
$ws.Range("L38").Value = 'Poste'
$ws.Range("M38").Value = -58.468425
$ws.Range("N38").Value = -34.54124

# Row 39
$ws.Range("A39").Value = '6020'
$ws.Range("B39").Value = '6/5/2025'
$ws.Range("C39").Value = 'RAVIGNANI, EMILIO, DR. 2036'
$ws.Range("D39").Value = '14'
$ws.Range("E39").Value = '807215465'
$ws.Range("F39").Value = 'NEW'
$ws.Range("G39").Value = 'Pendiente'
$ws.Range("H39").Value = 'Picada'
$ws.Range("I39").Value = '1'
$ws.Range("J39").Value = 'Cambio'
$ws.Range("K39").Value = 'Sin equipos'
$ws.Range("L39").Value = 'Pasante'
$ws.Range("M39").Value = -58.436298
$ws.Range("N39").Value = -34.578972

# Row 40
$ws.Range("A40").Value = '6144'
$ws.Range("B40").Value = '6/11/2025'
$ws.Range("C40").Value = 'TURIN 2990'
$ws.Range("D40").Value = '15'
$ws.Range("E40").Value = '807458282'
$ws.Range("F40").Value = 'NEW'
$ws.Range("G40").Value = 'Pendiente'
$ws.Range("H40").Value = 'Picada'
$ws.Range("I40").Value = '1'
$ws.Range("J40").Value = 'Cambio'
$ws.Range("K40").Value = 'Sin equipos'
$ws.Range("L40").Value = 'Pasante'
$ws.Range("M40").Value = -58.480925
$ws.Range("N40").Value = -34.585471

# Row 41
$ws.Range("A41").Value = '6143'
$ws.Range("B41").Value = '6/11/2025'
$ws.Range("C41").Value = 'SOLANO LOPEZ, F., MARISCAL 2845'
$ws.Range("D41").Value = '15'
$ws.Range("E41").Value = '807458296'
$ws.Range("F41").Value = 'NEW'
$ws.Range("G41").Value = 'Pendiente'
$ws.Range("H41").Value = 'Picada'
$ws.Range("I41").Value = '1'
$ws.Range("J41").Value = 'Cambio'
$ws.Range("K41").Value = 'Sin equipos'
$ws.Range("L41").Value = 'Pasante'
$ws.Range("M41").Value = -58.495071
$ws.Range("N41").Value = -34.593122

# Row 42
$ws.Range("A42").Value = '6141'
$ws.Range("B42").Value = '6/11/2025'
$ws.Range("C42").Value = 'EL PAMPERO 2618'
$ws.Range("D42").Value = '11'
$ws.Range("E42").Value = '807458310'
$ws.Range("F42").Value = 'NEW'
$ws.Range("G42").Value = 'Pendiente'
$ws.Range("H42").Value = 'Picada'
$ws.Range("I42").Value = '1'
$ws.Range("J42").Value = 'Cambio'
$ws.Range("K42").Value = 'Sin equipos'
$ws.Range("L42").Value = 'Pasante'
$ws.Range("M42").Value = -58.481942
$ws.Range("N42").Value = -34.602989

# Row 43
$ws.Range("A43").Value = '-478'
$ws.Range("B43").Value = '6/15/2025'
$ws.Range("C43").Value = 'Chivilcoy 4875'
$ws.Range("D43").Value = '11'
$ws.Range("E43").Value = '807508509'
$ws.Range("F43").Value = 'NEW'
$ws.Range("G43").Value = 'Pendiente'
$ws.Range("H43").Value = 'Poste podrido'
$ws.Range("I43").Value = '1'
$ws.Range("J43").Value = 'Cambio'
$ws.Range("K43").Value = 'Sin equipos'
$ws.Range("L43").Value = 'Poste'
$ws.Range("M43").Value = -58.517389
$ws.Range("N43").Value = -34.593541

# Row 44
$ws.Range("A44").Value = '6171'
$ws.Range("B44").Value = '6/18/2025'
$ws.Range("C44").Value = 'CABELLO 3486'
$ws.Range("D44").Value = '14'
$ws.Range("E44").Value = '807658640'
$ws.Range("F44").Value = 'NEW'
$ws.Range("G44").Value = 'Pendiente'
$ws.Range("H44").Value = 'Columna inclinada evaluar con inspector un corrimiento'
$ws.Range("I44").Value = '1'
$ws.Range("J44").Value = 'Cambio'
$ws.Range("K44").Value = 'Sin equipos'
$ws.Range("L44").Value = 'Pasante'
$ws.Range("M44").Value = -58.409579
$ws.Range("N44").Value = -34.581134

# Row 45
$ws.Range("A45").Value = '6230'
$ws.Range("B45").Value = '6/24/2025'
$ws.Range("C45").Value = 'Santa maria de oro	2722'
$ws.Range("D45").Value = '14'
$ws.Range("E45").Value = '807763066'
$ws.Range("F45").Value = 'NEW'
$ws.Range("G45").Value = 'Pendiente'
$ws.Range("H45").Value = 'Picada'
$ws.Range("I45").Value = '1'
$ws.Range("J45").Value = 'Cambio'
$ws.Range("K45").Value = 'Sin equipos'
$ws.Range("L45").Value = 'Pasante'
$ws.Range("M45").Value = -58.422315
$ws.Range("N45").Value = -34.576988

# Row 46
$ws.Range("A46").Value = '6231'
$ws.Range("B46").Value = '6/24/2025'
$ws.Range("C46").Value = 'Ciudad de la Paz 275'
$ws.Range("D46").Value = '14'
$ws.Range("E46").Value = '807763070'
$ws.Range("F46").Value = 'NEW'
$ws.Range("G46").Value = 'Pendiente'
$ws.Range("H46").Value = 'Aplomar o cambiar'
$ws.Range("I46").Value = '1'
$ws.Range("J46").Value = 'Cambio'
$ws.Range("K46").Value = 'Sin equipos'
$ws.Range("L46").Value = 'Terminal'
$ws.Range("M46").Value = -58.441049
$ws.Range("N46").Value = -34.574625

# Row 47
$ws.Range("A47").Value = '6233'
$ws.Range("B47").Value = '6/24/2025'
$ws.Range("C47").Value = 'HUERGO 389'
$ws.Range("D47").Value = '11'
$ws.Range("E47").Value = '807763078'
$ws.Range("F47").Value = 'NEW'
$ws.Range("G47").Value = 'Pendiente'
$ws.Range("H47").Value = 'Picada'
$ws.Range("I47").Value = '1'
$ws.Range("J47").Value = 'Cambio'
$ws.Range("K47").Value = 'Sin equipos'
$ws.Range("L47").Value = 'Pasante'
$ws.Range("M47").Value = -58.432722
$ws.Range("N47").Value = -34.572371

# Row 48
$ws.Range("A48").Value = '-492'
$ws.Range("B48").Value = '6/26/2025'
$ws.Range("C48").Value = 'Correa 3087'
$ws.Range("D48").Value = '12'
$ws.Range("E48").Value = '807789705'
$ws.Range("F48").Value = 'NEW'
$ws.Range("G48").Value = 'Pendiente'
$ws.Range("H48").Value = 'Ver con prioridad Esta agarrado por la red de cobre'
$ws.Range("I48").Value = '1'
$ws.Range("J48").Value = 'Aplomo'
$ws.Range("K48").Value = 'Sin equipos'
$ws.Range("L48").Value = 'Poste'
$ws.Range("M48").Value = -58.479742
$ws.Range("N48").Value = -34.546292

# Row 49
$ws.Range("A49").Value = '-495'
$ws.Range("B49").Value = '6/30/2025'
$ws.Range("C49").Value = 'Ricardo Balbin 3827'
$ws.Range("D49").Value = '12'
$ws.Range("E49").Value = '807846850 '
$ws.Range("F49").Value = 'NEW'
$ws.Range("G49").Value = 'Pendiente'
$ws.Range("H49").Value = 'Picada sin fotos pasa Pedro'
$ws.Range("I49").Value = '0'
$ws.Range("J49").Value = 'Cambio'
$ws.Range("K49").Value = 'Sin equipos'
$ws.Range("L49").Value = 'Pasante'
$ws.Range("M49").Value = -58.484375
$ws.Range("N49").Value = -34.554597

# Row 50
$ws.Range("A50").Value = '-496'
$ws.Range("B50").Value = '6/30/2025'
$ws.Range("C50").Value = 'Ricardo Balbin 3851'
$ws.Range("D50").Value = '12'
$ws.Range("E50").Value = '807846856'
$ws.Range("F50").Value = 'NEW'
$ws.Range("G50").Value = 'Pendiente'
$ws.Range("H50").Value = 'Picada sin fotos pasa Pedro'
$ws.Range("I50").Value = '0'
$ws.Range("J50").Value = 'Cambio'
$ws.Range("K50").Value = 'Sin equipos'
$ws.Range("L50").Value = 'Pasante'
$ws.Range("M50").Value = -58.484761
$ws.Range("N50").Value = -34.554643

# Row 51
$ws.Range("A51").Value = '6279'
$ws.Range("B51").Value = '6/30/2025'
$ws.Range("C51").Value = 'CORREA 3850'
$ws.Range("D51").Value = '12'
$ws.Range("E51").Value = '807847111'
$ws.Range("F51").Value = 'NEW'
$ws.Range("G51").Value = 'Pendiente'
$ws.Range("H51").Value = 'Picada'
$ws.Range("I51").Value = '1'
$ws.Range("J51").Value = 'Cambio'
$ws.Range("K51").Value = 'Sin equipos'
$ws.Range("L51").Value = 'Pasante'
$ws.Range("M51").Value = -58.486971
$ws.Range("N51").Value = -34.550135

# Row 52
$ws.Range("A52").Value = '6295'
$ws.Range("B52").Value = '6/30/2025'
$ws.Range("C52").Value = 'SOLER 6017'
$ws.Range("D52").Value = '14'
$ws.Range("E52").Value = '807851636'
$ws.Range("F52").Value = 'NEW'
$ws.Range("G52").Value = 'Pendiente'
$ws.Range("H52").Value = 'Picada'
$ws.Range("I52").Value = '1'
$ws.Range("J52").Value = 'Cambio'
$ws.Range("K52").Value = 'Sin equipos'
$ws.Range("L52").Value = 'Pasante'
$ws.Range("M52").Value = -58.436808
$ws.Range("N52").Value = -34.577464
